# Metadeconfound R for report
# Re-creates the author's edit: the metadata sheet was sorted by the
# "Batch" column (E), which re-orders the B2/B3 batches within the
# CV*/HA* block (rows 14-28 in the original layout); an AutoFilter was
# turned on for the header row; the selection was left on B23:B28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sort A1:E29 (header in row 1) ascending by column E ("Batch") ---
# xlYes = 1 (there is a header row), xlAscending = 1, xlSortNormal = 0
$sortRange = $ws.Range("A1:E29")
$sortKey   = $ws.Range("E1:E1")
$sortRange.Sort($sortKey, 1, $null, $null, 1, $null, 1, 1)

# --- Turn on AutoFilter for the header row ---
[void]$ws.Range("A1:E1").AutoFilter()

# Excel registers a hidden workbook-level name that backs every AutoFilter
# range; recreate it so the defined name matches what Excel itself writes.
$filterDb = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$E`$1")
$filterDb.Visible = $false

# --- Match the final selection recorded in the sheet view ---
[void]$ws.Range("B23:B28").Select()

# --- Best-effort: restore the saved window geometry (cosmetic only) ---
try {
    $win = $wb.Windows.Item(1)
    $win.Left   = 0
    $win.Top    = 0
    $win.Width  = 25600
    $win.Height = 16000
} catch {
    # Window geometry isn't always settable in headless hosts; ignore.
}

Write-Host "Sorted by Batch, enabled AutoFilter, restored selection."
